# Update LR-pair statistics per commit "Natmi following Dr Hou advice"
# Ligand-/Receptor-expressing cell counts increase from 1 to 3 for every row,
# which changes the derived expression/specificity/weight values accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 21.84976866666667
$ws.Cells.Item(2, 8).Value = 65.549306
$ws.Cells.Item(2, 9).Value = 0.05020018890879543
$ws.Cells.Item(2, 10).Value = 0.05020018890879543
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 112.513392
$ws.Cells.Item(2, 14).Value = 337.540176
$ws.Cells.Item(2, 15).Value = 0.3275312977368564
$ws.Cells.Item(2, 16).Value = 0.3275312977368564
$ws.Cells.Item(2, 17).Value = 2458.391587101984
$ws.Cells.Item(2, 18).Value = 22125.52428391785
$ws.Cells.Item(2, 19).Value = 0.01644213301993311
$ws.Cells.Item(2, 20).Value = 0.01644213301993311

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 21.84976866666667
$ws.Cells.Item(3, 8).Value = 65.549306
$ws.Cells.Item(3, 9).Value = 0.05020018890879543
$ws.Cells.Item(3, 10).Value = 0.05020018890879543
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 106.314466
$ws.Cells.Item(3, 14).Value = 318.943398
$ws.Cells.Item(3, 15).Value = 0.3094859589441663
$ws.Cells.Item(3, 16).Value = 0.3094859589441664
$ws.Cells.Item(3, 17).Value = 2322.946488020199
$ws.Cells.Item(3, 18).Value = 20906.51839218179
$ws.Cells.Item(3, 19).Value = 0.01553625360361685
$ws.Cells.Item(3, 20).Value = 0.01553625360361686

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 21.84976866666667
$ws.Cells.Item(4, 8).Value = 65.549306
$ws.Cells.Item(4, 9).Value = 0.05020018890879543
$ws.Cells.Item(4, 10).Value = 0.05020018890879543
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 124.6916553333333
$ws.Cells.Item(4, 14).Value = 374.074966
$ws.Cells.Item(4, 15).Value = 0.3629827433189773
$ws.Cells.Item(4, 16).Value = 0.3629827433189773
$ws.Cells.Item(4, 17).Value = 2724.483823697066
$ws.Cells.Item(4, 18).Value = 24520.3544132736
$ws.Cells.Item(4, 19).Value = 0.01822180228524546
$ws.Cells.Item(4, 20).Value = 0.01822180228524546

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 385.0524703333334
$ws.Cells.Item(5, 8).Value = 1155.157411
$ws.Cells.Item(5, 9).Value = 0.8846641374295412
$ws.Cells.Item(5, 10).Value = 0.8846641374295412
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 112.513392
$ws.Cells.Item(5, 14).Value = 337.540176
$ws.Cells.Item(5, 15).Value = 0.3275312977368564
$ws.Cells.Item(5, 16).Value = 0.3275312977368564
$ws.Cells.Item(5, 17).Value = 43323.55953518271
$ws.Cells.Item(5, 18).Value = 389912.0358166444
$ws.Cells.Item(5, 19).Value = 0.2897551929935543
$ws.Cells.Item(5, 20).Value = 0.2897551929935543

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 385.0524703333334
$ws.Cells.Item(6, 8).Value = 1155.157411
$ws.Cells.Item(6, 9).Value = 0.8846641374295412
$ws.Cells.Item(6, 10).Value = 0.8846641374295412
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 106.314466
$ws.Cells.Item(6, 14).Value = 318.943398
$ws.Cells.Item(6, 15).Value = 0.3094859589441663
$ws.Cells.Item(6, 16).Value = 0.3094859589441664
$ws.Cells.Item(6, 17).Value = 40936.64776546918
$ws.Cells.Item(6, 18).Value = 368429.8298892226
$ws.Cells.Item(6, 19).Value = 0.2737911289158953
$ws.Cells.Item(6, 20).Value = 0.2737911289158954

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 385.0524703333334
$ws.Cells.Item(7, 8).Value = 1155.157411
$ws.Cells.Item(7, 9).Value = 0.8846641374295412
$ws.Cells.Item(7, 10).Value = 0.8846641374295412
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 124.6916553333333
$ws.Cells.Item(7, 14).Value = 374.074966
$ws.Cells.Item(7, 15).Value = 0.3629827433189773
$ws.Cells.Item(7, 16).Value = 0.3629827433189773
$ws.Cells.Item(7, 17).Value = 48012.82991605256
$ws.Cells.Item(7, 18).Value = 432115.4692444731
$ws.Cells.Item(7, 19).Value = 0.3211178155200916
$ws.Cells.Item(7, 20).Value = 0.3211178155200917

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 28.350479
$ws.Cells.Item(8, 8).Value = 85.05143699999999
$ws.Cells.Item(8, 9).Value = 0.06513567366166337
$ws.Cells.Item(8, 10).Value = 0.06513567366166337
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 112.513392
$ws.Cells.Item(8, 14).Value = 337.540176
$ws.Cells.Item(8, 15).Value = 0.3275312977368564
$ws.Cells.Item(8, 16).Value = 0.3275312977368564
$ws.Cells.Item(8, 17).Value = 3189.808557114768
$ws.Cells.Item(8, 18).Value = 28708.27701403291
$ws.Cells.Item(8, 19).Value = 0.02133397172336898
$ws.Cells.Item(8, 20).Value = 0.02133397172336898

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 28.350479
$ws.Cells.Item(9, 8).Value = 85.05143699999999
$ws.Cells.Item(9, 9).Value = 0.06513567366166337
$ws.Cells.Item(9, 10).Value = 0.06513567366166337
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 106.314466
$ws.Cells.Item(9, 14).Value = 318.943398
$ws.Cells.Item(9, 15).Value = 0.3094859589441663
$ws.Cells.Item(9, 16).Value = 0.3094859589441664
$ws.Cells.Item(9, 17).Value = 3014.066035729214
$ws.Cells.Item(9, 18).Value = 27126.59432156292
$ws.Cells.Item(9, 19).Value = 0.02015857642465416
$ws.Cells.Item(9, 20).Value = 0.02015857642465417

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 28.350479
$ws.Cells.Item(10, 8).Value = 85.05143699999999
$ws.Cells.Item(10, 9).Value = 0.06513567366166337
$ws.Cells.Item(10, 10).Value = 0.06513567366166337
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 124.6916553333333
$ws.Cells.Item(10, 14).Value = 374.074966
$ws.Cells.Item(10, 15).Value = 0.3629827433189773
$ws.Cells.Item(10, 16).Value = 0.3629827433189773
$ws.Cells.Item(10, 17).Value = 3535.068156002905
$ws.Cells.Item(10, 18).Value = 31815.61340402614
$ws.Cells.Item(10, 19).Value = 0.02364312551364023
$ws.Cells.Item(10, 20).Value = 0.02364312551364023
